# Add a "riotIdGameName" username column before the existing "game_minutes"
# column (new column B), shifting every other stat column one to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; this shifts B:M -> C:N (and values/styles
# move with it).
$ws.Range("B1").EntireColumn.Insert()

# New header cell should look like the rest of row 1 (bold/bordered style).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "riotIdGameName"

# Player usernames (row -> riotIdGameName), in row order.
$names = @{
    2  = "Chaffles"
    3  = "IamClone"
    4  = "백지 소라"
    5  = "Ganjegreen"
    6  = "Drewsph"
    7  = "MAA DuckJugs"
    8  = "Empro"
    9  = "redzawsome"
    10 = "Oasis RSexy"
    11 = "Doretha728"
    12 = "HPZ Tea Jay"
    13 = "JustaWittleGuy"
    14 = "CrimsonYoni"
    15 = "agateo"
    16 = "Abraxo"
    17 = "Oasis psychotikk"
    18 = "Bug"
    19 = "SimonLowell"
    20 = "bezzaboyo"
    21 = "Oasis Jags"
}

foreach ($row in $names.Keys) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $names[$row]
    # Data rows in this sheet carry no explicit style (unlike column A / the
    # header row), so strip the style the column-insert inherited from A.
    $cell.ClearFormats()
}
